$sheet1New = @(
    @(14, 44005, 28631, 691, 16006),
    @(15, 44012, 33387, 754, 17904),
    @(16, 44019, 39588, 829, 20056),
    @(17, 44026, 47671, 929, 23459)
)

$sheet2Week1 = @(
    @(418, 4, 831.11, 133),
    @(419, 5, 289.31, 10),
    @(420, 6, 27.72, 1),
    @(421, 7, 305.52999999999997, 8),
    @(422, 8, 49.65, 2),
    @(423, 9, 300.42, 87),
    @(424, 10, 130.91, 1),
    @(425, 11, 49.19, 0),
    @(426, 12, 313.95, 14),
    @(427, 13, 147.19999999999999, 1),
    @(428, 14, 230.61, 2),
    @(429, 15, 317.32, 10),
    @(430, 16, 319.24, 37),
    @(431, 17, 159.46, 4),
    @(432, 18, 51.19, 3),
    @(433, 19, 162.87, 0),
    @(434, 20, 104.85, 8),
    @(435, 21, 157.25, 37),
    @(436, 22, 293.01, 18),
    @(437, 23, 62.41, 1),
    @(438, 24, 301.83999999999997, 31),
    @(439, 25, 158.75, 8),
    @(440, 26, 97.72, 14),
    @(441, 27, 219.11, 10),
    @(442, 28, 229.15, 104),
    @(443, 29, 144.94, 4),
    @(444, 30, 132.97999999999999, 13),
    @(445, 31, 245.01, 7),
    @(446, 32, 34.03, 7),
    @(447, 33, 38.479999999999997, 2),
    @(448, 34, 136.02000000000001, 3),
    @(449, 35, 339.23, 174)
)

$sheet2Week2 = @(
    @(450, 4, 1035.96, 149),
    @(451, 5, 340.15, 11),
    @(452, 6, 43.56, 2),
    @(453, 7, 342.01, 8),
    @(454, 8, 75.23, 2),
    @(455, 9, 320.49, 88),
    @(456, 10, 171.92, 1),
    @(457, 11, 66.3, 0),
    @(458, 12, 349.39, 15),
    @(459, 13, 219.09, 1),
    @(460, 14, 288.62, 2),
    @(461, 15, 362.81, 10),
    @(462, 16, 362.52, 39),
    @(463, 17, 186.4, 5),
    @(464, 18, 59.72, 3),
    @(465, 19, 274.31, 0),
    @(466, 20, 127.13, 9),
    @(467, 21, 184.86, 41),
    @(468, 22, 312.54000000000002, 18),
    @(469, 23, 77.569999999999993, 1),
    @(470, 24, 339.8, 38),
    @(471, 25, 211.51, 8),
    @(472, 26, 124.6, 14),
    @(473, 27, 255.96, 10),
    @(474, 28, 293.35000000000002, 110),
    @(475, 29, 151.91999999999999, 4),
    @(476, 30, 139.77000000000001, 13),
    @(477, 31, 270.83, 12),
    @(478, 32, 39.78, 7),
    @(479, 33, 51.3, 2),
    @(480, 34, 185.64, 5),
    @(481, 35, 390.43, 201)
)

$sheet2Week3 = @(
    @(482, 4, 1257.3, 169),
    @(483, 5, 404.04, 11),
    @(484, 6, 45.54, 3),
    @(485, 7, 357.34, 9),
    @(486, 8, 102.31, 2),
    @(487, 9, 361.64, 89),
    @(488, 10, 209.77, 1),
    @(489, 11, 112.29, 0),
    @(490, 12, 411.93, 22),
    @(491, 13, 220.8, 1),
    @(492, 14, 398.13, 2),
    @(493, 15, 420.77, 10),
    @(494, 16, 414.79, 43),
    @(495, 17, 236.01, 5),
    @(496, 18, 73.37, 5),
    @(497, 19, 320.02999999999997, 0),
    @(498, 20, 166.13, 9),
    @(499, 21, 217.87, 45),
    @(500, 22, 333.16, 18),
    @(501, 23, 171.18, 1),
    @(502, 24, 366.47, 45),
    @(503, 25, 289.98, 9),
    @(504, 26, 169.86, 14),
    @(505, 27, 353.54, 10),
    @(506, 28, 379.27, 118),
    @(507, 29, 195.58, 4),
    @(508, 30, 142.6, 13),
    @(509, 31, 298.77, 13),
    @(510, 32, 52.87, 7),
    @(511, 33, 66.459999999999994, 3),
    @(512, 34, 191.16, 5),
    @(513, 35, 444.25, 243)
)

$provinceNames = @{
    4 = "Distrito Nacional"
    5 = "Azua"
    6 = "Baoruco"
    7 = "Barahona"
    8 = "Dajabon"
    9 = "Duarte"
    10 = "Elias Pina"
    11 = "El Seibo"
    12 = "Espaillat"
    13 = "Independencia"
    14 = "La Altagracia"
    15 = "La Romana"
    16 = "La Vega"
    17 = "Maria Trinidad Sanchez"
    18 = "Monte Cristi"
    19 = "Pedernales"
    20 = "Peravia"
    21 = "Puerto Plata"
    22 = "Hermanas Mirabal"
    23 = "Samana"
    24 = "San Cristobal"
    25 = "San Juan"
    26 = "San Pedro de Macoris"
    27 = "Sanchez Ramirez"
    28 = "Santiago"
    29 = "Santiago Rodriguez"
    30 = "Valverde"
    31 = "Monsenor Nouel"
    32 = "Monte Plata"
    33 = "Hato Mayor"
    34 = "San Jose de Ocoa"
    35 = "Santo Domingo"
}

$ageNames = @{
    37 = "<1"
    38 = "20-29"
    39 = "30-39"
    40 = "40-49"
    41 = "50-59"
    42 = "60+"
    43 = "1-4"
    44 = "5-9"
    45 = "10-19"
}

$sheet3New = @(
    @(65, 44012, 37, 0, 0.4, 1.25),
    @(66, 44012, 43, 1, 1.5, 0.45),
    @(67, 44012, 44, 1, 1.9, 0),
    @(68, 44012, 45, 1, 4.8, 2.1),
    @(69, 44012, 38, 0, 19.149999999999999, 6.6),
    @(70, 44012, 39, 0, 22.45, 5.8),
    @(71, 44012, 40, 0, 17.5, 8.5),
    @(72, 44012, 41, 0, 13.85, 15.9),
    @(73, 44012, 42, 0, 15.05, 63.55),
    @(74, 44019, 37, 0, 0.45, 1.1000000000000001),
    @(75, 44019, 43, 1, 1.55, 0.4),
    @(76, 44019, 44, 1, 1.9, 0.2),
    @(77, 44019, 45, 1, 4.75, 0.95),
    @(78, 44019, 38, 0, 19.399999999999999, 2.95),
    @(79, 44019, 39, 0, 22.6, 5.7),
    @(80, 44019, 40, 0, 17.45, 8.5500000000000007),
    @(81, 44019, 41, 0, 13.7, 16.2),
    @(82, 44019, 42, 0, 14.8, 63.8),
    @(83, 44026, 37, 0, 0.45, 0.95),
    @(84, 44026, 43, 1, 1.65, 0.35),
    @(85, 44026, 44, 1, 1.9, 0.15),
    @(86, 44026, 45, 1, 4.7, 0.8),
    @(87, 44026, 38, 0, 19.45, 2.8),
    @(88, 44026, 39, 0, 22.8, 5.6),
    @(89, 44026, 40, 0, 17.5, 9.35),
    @(90, 44026, 41, 0, 13.55, 15.1),
    @(91, 44026, 42, 0, 14.7, 64.75)
)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Fallecido_Recuperado" (sheet1): append weekly totals rows 14-17
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Fallecido_Recuperado")

# Copy the date style (s="2") from the last existing data row down onto the
# new rows so the new date cells reuse the existing numFmt style instead of
# creating a new one.
$ws1.Range("A13").Copy()
$ws1.Range("A14:A17").PasteSpecial(-4122)

foreach ($row in $sheet1New) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
}

$ws1.Range("D18").Select()

# ---------------------------------------------------------------------------
# Sheet "Provincias_Semanal" (sheet2): fix the mislabeled 6/24 week (it was
# stored as literal text "24/6/2020" instead of a real date) and append the
# two new weeks of province data.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Provincias_Semanal")

# Reuse the date-style format from the prior week's A column across rows
# 418-513 (the fixed week + the two newly appended weeks) in one shot.
$ws2.Range("A417").Copy()
$ws2.Range("A418:A513").PasteSpecial(-4122)

foreach ($row in $sheet2Week1) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = 44012
    $ws2.Cells.Item($r, 2).Value = $provinceNames[[int]$row[1]]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
}

foreach ($row in $sheet2Week2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = 44019
    $ws2.Cells.Item($r, 2).Value = $provinceNames[[int]$row[1]]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
}

foreach ($row in $sheet2Week3) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = 44026
    $ws2.Cells.Item($r, 2).Value = $provinceNames[[int]$row[1]]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
}

$ws2.Activate()
$ws2.Range("D514").Select()
$excel.ActiveWindow.ScrollRow = 486

# ---------------------------------------------------------------------------
# Sheet "Por_Edad" (sheet3): append the three new weeks of age-bracket data.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Por_Edad")

$ws3.Range("A64").Copy()
$ws3.Range("A65:A91").PasteSpecial(-4122)

foreach ($row in $sheet3New) {
    $r = $row[0]
    $date = $row[1]
    $ageIdx = [int]$row[2]
    $style = [int]$row[3]
    $c = $row[4]
    $d = $row[5]

    $ws3.Cells.Item($r, 1).Value = $date

    if ($style -eq 1) {
        $ws3.Range("B57").Copy()
        $ws3.Range("B$r").PasteSpecial(-4122)
        $ws3.Cells.Item($r, 2).Value = "'" + $ageNames[$ageIdx]
    } else {
        $ws3.Cells.Item($r, 2).Value = $ageNames[$ageIdx]
    }

    $ws3.Cells.Item($r, 3).Value = $c
    $ws3.Cells.Item($r, 4).Value = $d
}

$ws3.Activate()
$ws3.Range("D92").Select()
$excel.ActiveWindow.ScrollRow = 67

$ws1.Activate()
